$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Constants
# ---------------------------------------------------------------------------
$xlNone             = -4142   # xlLineStyleNone
$xlContinuous       = 1       # xlContinuous
$xlEdgeLeft         = 7
$xlEdgeTop          = 8
$xlEdgeBottom       = 9
$xlEdgeRight        = 10
$xlPasteFormats     = -4122   # xlPasteFormats

$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# ---------------------------------------------------------------------------
# Build two throw-away "template" cells holding the exact two new border
# styles needed (top+bottom only, and top+bottom+right), then stamp copies
# of those formats onto every target cell. Doing it this way (one finished
# format built once, then replicated with PasteSpecial) keeps the workbook's
# style table tidy - exactly two new cellXfs get added overall.
# ---------------------------------------------------------------------------

# Template 1: top + bottom border only
$tmpl1 = $ws1.Range("Z1")
$tmpl1.ClearFormats()
$tmpl1.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$tmpl1.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
$tmpl1.Copy()

$ws1.Range("C1").PasteSpecial($xlPasteFormats)
$ws2.Range("C1").PasteSpecial($xlPasteFormats)
$ws2.Range("F1").PasteSpecial($xlPasteFormats)

$tmpl1.ClearFormats()
$tmpl1.Clear()

# Template 2: top + bottom + right border
$tmpl2 = $ws1.Range("Z2")
$tmpl2.ClearFormats()
$tmpl2.Borders.LineStyle = $xlContinuous
$tmpl2.Borders.Item($xlEdgeLeft).LineStyle = $xlNone
$tmpl2.Copy()

$ws1.Range("D1").PasteSpecial($xlPasteFormats)
$ws2.Range("D1").PasteSpecial($xlPasteFormats)
$ws2.Range("G1").PasteSpecial($xlPasteFormats)

$tmpl2.ClearFormats()
$tmpl2.Clear()

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Text updates: "fedcore" -> "approach" (anonymisation)
# ---------------------------------------------------------------------------
$ws1.Range("C2").Value = "approach"

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# ---------------------------------------------------------------------------
# Drop the stray empty inline-string cell at G5 on computational_comparison
# ---------------------------------------------------------------------------
$ws2.Range("G5").ClearContents()
